# Fix property_category values:
#  - "建物" (building) sheet: column I (property_category) rows 2-7 should be "building" instead of "land"
#  - "汽車" (car) sheet: column H (property_category) rows 2-3 should be "car" instead of "land"

$wb = $excel.ActiveWorkbook

$wsBuilding = $wb.Worksheets.Item("建物")
for ($r = 2; $r -le 7; $r++) {
    $wsBuilding.Cells.Item($r, 9).Value = "building"
}

$wsCar = $wb.Worksheets.Item("汽車")
for ($r = 2; $r -le 3; $r++) {
    $wsCar.Cells.Item($r, 8).Value = "car"
}
